$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 5161.07
$ws.Range("E2").Value = -5161.07

$ws.Range("D4").Value = 10275.72
$ws.Range("E4").Value = 3447.620000000001
$ws.Range("F4").Value = 0.7487769012499872
